$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Turni")

# --- Row 8 ---
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 5783861406
$ws.Cells.Item(8,3).Value = "Oooo"
$ws.Cells.Item(8,4).Value = "Vvvvv"
$ws.Cells.Item(8,5).Value = 18
$ws.Cells.Item(8,6).Value = "CERCA 1"

# Column G looks like a bare date ("2025-12-17") - Excel would silently
# convert that to a date serial number on assignment, so force a text
# number format, assign it, then clear the format back off again so the
# stored value stays the literal string without leaving a style behind.
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "2025-12-17"
$ws.Cells.Item(8,7).ClearFormats()

$ws.Cells.Item(8,8).Value = "2025-12-17 13:30:10"
$ws.Cells.Item(8,9).Value = "2025-12-17 13:30:41"
$ws.Cells.Item(8,10).Value = 0.01
$ws.Cells.Item(8,11).Value = "d:\MO.VE_project_CA__Pt\Pulizie_BOT_MOVE\archivio\video\2025\12\17\CERCA_1\Oooo_Vvvvv_ingresso_13-30.mp4"
$ws.Cells.Item(8,12).Value = "BAACAgQAAxkBAAIBrWlColIVKlGEod-hBK9cr3CQ_6nWAAJ4GAAC6vYYUi98IQuw9AXHNgQ"
$ws.Cells.Item(8,13).Value = "d:\MO.VE_project_CA__Pt\Pulizie_BOT_MOVE\archivio\video\2025\12\17\CERCA_1\Oooo_Vvvvv_uscita_13-30.mp4"
$ws.Cells.Item(8,14).Value = "BAACAgQAAxkBAAIBuGlConGwbQwY97zxMHC74d3bd6MOAAJ5GAAC6vYYUgsIX-xdL-UVNgQ"
$ws.Cells.Item(8,15).Value = "completato"

# --- Row 9 ---
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 5783861406
$ws.Cells.Item(9,3).Value = "Oooo"
$ws.Cells.Item(9,4).Value = "Vvvvv"
$ws.Cells.Item(9,5).Value = 65
$ws.Cells.Item(9,6).Value = "Rainusso"

$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "2025-12-17"
$ws.Cells.Item(9,7).ClearFormats()

$ws.Cells.Item(9,8).Value = "2025-12-17 15:28:58"
$ws.Cells.Item(9,9).Value = "2025-12-17 15:57:35"
$ws.Cells.Item(9,10).Value = 0.48
$ws.Cells.Item(9,11).Value = "D:\MO.VE_project_CA__Pt\Pulizie_BOT_MOVE\archivio\video\2025\12\17\Rainusso\Oooo_Vvvvv_ingresso_15-28.mp4"
$ws.Cells.Item(9,12).Value = "BAACAgQAAxkBAAICWWlCvioMOEb0mqcXsuYXExLEAZCuAALKGAAC6vYYUlwK0P7tglcQNgQ"
$ws.Cells.Item(9,13).Value = "d:\MO.VE_project_CA__Pt\Pulizie_BOT_MOVE\archivio\video\2025\12\17\Rainusso\Oooo_Vvvvv_uscita_15-57.mp4"
$ws.Cells.Item(9,14).Value = "BAACAgQAAxkBAAICa2lCxN-ap43h5H2fk7Ca0SH-qgFCAALZGAAC6vYYUhRaPDP1iRGHNgQ"
$ws.Cells.Item(9,15).Value = "completato"
